$wb = $excel.ActiveWorkbook

# Update "想去人数" (number of people interested) values for two rows that
# appear on both the "展览" sheet and the "全部类型" sheet.
$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("F4").Value = 26
    $ws.Range("F5").Value = 2326
}
